$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = "'34.788.62"
$ws.Range("E2").Value = '  -1.99%  '

# Row 3
$ws.Range("D3").Formula = "'1.872.44"
$ws.Range("E3").Value = '  -2.80%  '

# Row 4
$ws.Range("E4").Value = '  -0.72%  '

# Row 5
$ws.Range("D5").Formula = "'245.91"
$ws.Range("E5").Value = '  -3.40%  '

# Row 6
$ws.Range("D6").Formula = "'0.680"
$ws.Range("E6").Value = '  -7.49%  '

# Row 8
$ws.Range("D8").Formula = "'41.71"
$ws.Range("E8").Value = '  +1.85%  '

# Row 9
$ws.Range("D9").Formula = "'0.344"
$ws.Range("E9").Value = '  -3.96%  '

# Row 10
$ws.Range("D10").Formula = "'51.01"
$ws.Range("E10").Value = '  -3.37%  '

# Row 11
$ws.Range("D11").Formula = "'0.0732"
$ws.Range("E11").Value = '  -2.97%  '

# Row 12
$ws.Range("D12").Formula = "'0.0969"
$ws.Range("E12").Value = '  -3.23%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Formula = "'12.84"
$ws.Range("E13").Value = '  -0.38%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Formula = "'2.145.62"
$ws.Range("E14").Value = '  -2.73%  '

# Row 15
$ws.Range("D15").Formula = "'0.709"
$ws.Range("E15").Value = '  -1.95%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Formula = "'1.876.87"
$ws.Range("E16").Value = '  -2.83%  '

# Row 17
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").Formula = "'4.83"
$ws.Range("E17").Value = '  -1.98%  '

# Row 18
$ws.Range("D18").Formula = "'34.826.79"
$ws.Range("E18").Value = '  -1.92%  '

# Row 19
$ws.Range("D19").Formula = "'72.34"
$ws.Range("E19").Value = '  -3.15%  '

# Row 20
$ws.Range("D20").Formula = "'" + '0.0' + [char]0x2083 + '0813'
$ws.Range("E20").Value = '  -3.13%  '

# Row 21
$ws.Range("D21").Formula = "'243.02"
$ws.Range("E21").Value = '  -0.80%  '

# Row 22
$ws.Range("D22").Formula = "'12.58"
$ws.Range("E22").Value = '  -3.56%  '

# Row 23
$ws.Range("D23").Formula = "'4.88"
$ws.Range("E23").Value = '  -4.66%  '

# Row 24
$ws.Range("E24").Value = '  -0.77%  '

# Row 25
$ws.Range("D25").Formula = "'2.49"
$ws.Range("E25").Value = '  +5.88%  '

# Row 26
$ws.Range("D26").Formula = "'2.19"
$ws.Range("E26").Value = '  -10.61%  '

# Row 27
$ws.Range("D27").Formula = "'164.25"
$ws.Range("E27").Value = '  -2.03%  '

# Row 28
$ws.Range("D28").Formula = "'8.32"
$ws.Range("E28").Value = '  -3.84%  '

# Row 29
$ws.Range("D29").Formula = "'18.10"
$ws.Range("E29").Value = '  -4.03%  '

# Row 30
$ws.Range("E30").Value = '  -6.97%  '

# Row 31
$ws.Range("D31").Formula = "'4.128.38"
$ws.Range("E31").Value = '  -0.07%  '

# Row 32
$ws.Range("D32").Formula = "'1.68"
$ws.Range("E32").Value = '  +0.99%  '

# Row 33
$ws.Range("D33").Formula = "'4.18"
$ws.Range("E33").Value = '  -4.40%  '

# Row 34
$ws.Range("D34").Formula = "'0.0574"
$ws.Range("E34").Value = '  -1.67%  '

# Row 35
$ws.Range("B35").Value = 'BinanceUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D35").Formula = "'1.00"
$ws.Range("E35").Value = '  -0.84%  '

# Row 36
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").Formula = "'4.11"
$ws.Range("E36").Value = '  -3.37%  '

# Row 37
$ws.Range("D37").Formula = "'0.822"
$ws.Range("E37").Value = '  -10.52%  '

# Row 38
$ws.Range("D38").Formula = "'1.97"
$ws.Range("E38").Value = '  -2.89%  '

# Row 39
$ws.Range("E39").Value = '  -22.90%  '

# Row 40
$ws.Range("D40").Formula = "'97.06"
$ws.Range("E40").Value = '  -0.57%  '

# Row 41
$ws.Range("D41").Formula = "'16.89"
$ws.Range("E41").Value = '  -3.23%  '

# Row 42
$ws.Range("D42").Formula = "'0.0665"
$ws.Range("E42").Value = '  +2.23%  '

# Row 43
$ws.Range("D43").Formula = "'0.0210"
$ws.Range("E43").Value = '  -0.55%  '

# Row 44
$ws.Range("E44").Value = '  -5.60%  '

# Row 45
$ws.Range("D45").Formula = "'1.280.42"
$ws.Range("E45").Value = '  -5.23%  '

# Row 46
$ws.Range("D46").Formula = "'0.0817"
$ws.Range("E46").Value = '  +10.97%  '

# Row 47
$ws.Range("D47").Formula = "'2.29"
$ws.Range("E47").Value = '  -7.49%  '

# Row 48
$ws.Range("E48").Value = '  -1.01%  '

# Row 49
$ws.Range("E49").Value = '  -1.98%  '

# Row 50
$ws.Range("D50").Formula = "'6.36"
$ws.Range("E50").Value = '  -7.13%  '

# Row 51
$ws.Range("D51").Formula = "'11.76"
$ws.Range("E51").Value = '  -1.68%  '
